$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed contact list only has 31 records (rows 2-32); the old
# sheet had 47 (rows 2-48), so drop the now-unused trailing rows first.
# EntireRow.Delete() shifts the rows below up and keeps the dimension in sync.
$ws.Range("A33:C48").EntireRow.Delete()

# Phone numbers that are purely numeric must stay text (matching the source
# data) instead of being auto-converted to numbers, so force column C to a
# text format before assigning those values ("N/A" rows are already text).
$ws.Range("C2:C17").NumberFormat = "@"
$ws.Range("C19:C31").NumberFormat = "@"

# Replace rows 2-32 with the new Belgium contact data.
$ws.Range("A2").Value = 'IGOR JOVANOVIC'
$ws.Range("B2").Value = 'Rue du Commerce 51, 4100 Seraing'
$ws.Range("C2").Value = '3243396154'
$ws.Range("A3").Value = 'Verica Jovanovic'
$ws.Range("B3").Value = 'Paardenmarktstraat 129, 3080 Tervuren'
$ws.Range("C3").Value = '3227680286'
$ws.Range("A4").Value = 'Jovanovic Dragica'
$ws.Range("B4").Value = 'Rue Général Henry 134, 1040 Etterbeek'
$ws.Range("C4").Value = '3223100295'
$ws.Range("A5").Value = 'Jovanovic Draga'
$ws.Range("B5").Value = 'Akkerstraat 23, 3680 Neeroeteren (Maaseik)'
$ws.Range("C5").Value = '3289215951'
$ws.Range("A6").Value = 'Nikola Jovanovic'
$ws.Range("B6").Value = 'Drève des Equipages 16, 1170 Watermael-Boitsfort'
$ws.Range("C6").Value = '3226602204'
$ws.Range("A7").Value = 'Joey Jovanovic'
$ws.Range("B7").Value = 'Rue des Etangs Noirs 89, 1080 Molenbeek-Saint-Jean'
$ws.Range("C7").Value = '3226404787'
$ws.Range("A8").Value = 'DRAGANA JOVANOVIC'
$ws.Range("B8").Value = 'Rue Péterson 7, 5580 Jemelle (Rochefort)'
$ws.Range("C8").Value = '3284730307'
$ws.Range("A9").Value = 'SVETLANA JOVANOVIC'
$ws.Range("B9").Value = 'Rue du Warmonceau 33, 6061 Montignies-sur-Sambre (Charleroi)'
$ws.Range("C9").Value = '3271120443'
$ws.Range("A10").Value = 'SIBIN JOVANOVIC'
$ws.Range("B10").Value = 'Kapellenboslaan 35, 2830 Willebroek'
$ws.Range("C10").Value = '3234374396'
$ws.Range("A11").Value = 'DAVID JOVANOVIC'
$ws.Range("B11").Value = 'Rue de Montegnée 2, 4101 Jemeppe-sur-Meuse (Seraing)'
$ws.Range("C11").Value = '3243911905'
$ws.Range("A12").Value = 'Jika Jovanovic'
$ws.Range("B12").Value = 'Petite Rue du Moulin 26, 1070 Anderlecht'
$ws.Range("C12").Value = '3226444124'
$ws.Range("A13").Value = 'RENATE JOVANOVIC'
$ws.Range("B13").Value = 'Domaine des Maraîchers 68, 7390 Wasmuel (Quaregnon)'
$ws.Range("C13").Value = '3265641298'
$ws.Range("A14").Value = 'ANITA JOVANOVIC'
$ws.Range("B14").Value = 'KERKENDIJK 19, 2000 ANTWERPEN'
$ws.Range("C14").Value = '3234370683'
$ws.Range("A15").Value = 'Tycha Jovanovic'
$ws.Range("B15").Value = 'Rue du Fond 8, 4920 Sougné-Remouchamps (Aywaille)'
$ws.Range("C15").Value = '3242431388'
$ws.Range("A16").Value = 'MILAN JOVANOVIC'
$ws.Range("B16").Value = 'Rue Jean Noté 42, 1070 Anderlecht'
$ws.Range("C16").Value = '3225234241'
$ws.Range("A17").Value = 'LAORA JOVANOVIC'
$ws.Range("B17").Value = 'Rue Arthur Decoux 54, 6020 Dampremy (Charleroi)'
$ws.Range("C17").Value = '3271575700'
$ws.Range("A18").Value = 'Jovanovic Branko'
$ws.Range("B18").Value = 'Avenue des Sorbiers 11, 1330 Rixensart'
$ws.Range("C18").Value = 'N/A'
$ws.Range("A19").Value = 'DZULISTANA JOVANOVIC'
$ws.Range("B19").Value = 'Rue Edmond Leburton 14, 6200 Châtelineau (Châtelet)'
$ws.Range("C19").Value = '3271972091'
$ws.Range("A20").Value = 'Dragana Jovanovic'
$ws.Range("B20").Value = 'Rue des Alouettes 32, 7100 La Louvière'
$ws.Range("C20").Value = '3264222958'
$ws.Range("A21").Value = 'DARKO JOVANOVIC'
$ws.Range("B21").Value = 'KESSELDAALLAAN 8, 3010 KESSEL-LO'
$ws.Range("C21").Value = '3216427363'
$ws.Range("A22").Value = 'Draga Jovanovic'
$ws.Range("B22").Value = 'Dascottelei 55/401, 2100 Deurne (Anvers)'
$ws.Range("C22").Value = '3233211267'
$ws.Range("A23").Value = 'SANJA JOVANOVIC'
$ws.Range("B23").Value = 'Kampioenstraat 26, 2020 Anvers'
$ws.Range("C23").Value = '3232764929'
$ws.Range("A24").Value = 'Srdan Jovanovic'
$ws.Range("B24").Value = 'Daliastraat 12, 2580 Putte'
$ws.Range("C24").Value = '3215760287'
$ws.Range("A25").Value = 'DOBRICA JOVANOVIC'
$ws.Range("B25").Value = 'Rue des Pâquerettes 52, 7160 Chapelle-lez-Herlaimont'
$ws.Range("C25").Value = '3264382648'
$ws.Range("A26").Value = 'DARKO JOVANOVIC'
$ws.Range("B26").Value = 'Kesseldallaan 8, 3010 Kessel-Lo (Louvain)'
$ws.Range("C26").Value = '3216750744'
$ws.Range("A27").Value = 'Marija Jovanovic'
$ws.Range("B27").Value = 'Oranjestraat 31 R, 2060 Anvers'
$ws.Range("C27").Value = '3232319317'
$ws.Range("A28").Value = 'Isnija Jovanovic'
$ws.Range("B28").Value = 'Rue de l''Ourthe 31, 1080 Molenbeek-Saint-Jean'
$ws.Range("C28").Value = '3226443721'
$ws.Range("A29").Value = 'Jelica Jovanovic'
$ws.Range("B29").Value = 'Laurian-Moris-Straße 17, 4780 Saint-Vith'
$ws.Range("C29").Value = '3280685676'
$ws.Range("A30").Value = 'Samantha Jovanovic'
$ws.Range("B30").Value = 'Boulevard Joseph Tirou 10 00O7, 6000 Charleroi'
$ws.Range("C30").Value = '3271819213'
$ws.Range("A31").Value = 'Dusko Jovanovic'
$ws.Range("B31").Value = 'Rue Willems 45/5, 1210 Saint-Josse-ten-Noode'
$ws.Range("C31").Value = '3223301843'
$ws.Range("A32").Value = 'Milorad Jovanovic'
$ws.Range("B32").Value = 'Rue de l''Harmonie 7 b007, 1000 Bruxelles'
$ws.Range("C32").Value = 'N/A'
